$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 737, shifting existing rows 737:778 down to 738:779.
$ws.Rows.Item(737).Insert()

# Populate the newly inserted row with the new data point for 2026/02/01.
# Force column A to text first so the date-like string isn't auto-converted
# to a date serial (matches the rest of the sheet, which stores dates as
# plain text), then restore the default "Normal" style so no stray
# formatting is left behind.
$ws.Cells.Item(737, 1).NumberFormat = "@"
$ws.Cells.Item(737, 1).Value = "2026/02/01"
$ws.Cells.Item(737, 1).Style = "Normal"
$ws.Cells.Item(737, 2).Value = "日"
$ws.Cells.Item(737, 3).Value = 19
$ws.Cells.Item(737, 4).Value = 22
